# Sales report fix-up:
#  - Order #FBZ-L1V4H3K (Holistic Fitness ... Shipped) was recorded after
#    #FBZ-L1RGCSK (MEDIX ... Pending) even though its order date is earlier;
#    same mix-up for Joyal K's two orders further down. Put them back in the
#    right order.
#  - A missing order (#FBZ-328BOQA, Achu K) is added as a new row, and the
#    Summary totals at the bottom are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the cell to be stored as text even when the value looks numeric
    # (e.g. "829", "4,599", "0"), matching how the sheet already stores the
    # Amount/Discount columns. Reset back to the default "Normal" style
    # afterwards so no lingering number-format is left on the cell.
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# ---- Fix the swapped order rows -----------------------------------------

# Row 2 <-> Row 3 (Joyal Kuriakose's two orders were in the wrong order)
Set-TextCell "B2" "#FBZ-L1V4H3K"
$ws.Range("C2").Value = "Joyal Kuriakose"
Set-TextCell "D2" "2025-02-04"
$ws.Range("E2").Value = "Holistic Fitness 6kg Rubber Slam Ball"
Set-TextCell "G2" "4,599"
Set-TextCell "H2" "0"
$ws.Range("I2").Value = "netbanking"
$ws.Range("J2").Value = "Shipped"

Set-TextCell "B3" "#FBZ-L1RGCSK"
$ws.Range("C3").Value = "Joyal Kuriakose"
Set-TextCell "D3" "2025-01-04"
$ws.Range("E3").Value = "MEDIX Soft Medicine Ball (2), Rubber for Adults"
Set-TextCell "G3" "829"
Set-TextCell "H3" "0"
$ws.Range("I3").Value = "netbanking"
$ws.Range("J3").Value = "Pending"

# Row 4 <-> Row 5 (Joyal K's two orders were in the wrong order)
Set-TextCell "B4" "#FBZ-6C73256"
$ws.Range("C4").Value = "Joyal K"
Set-TextCell "D4" "2025-04-04"
$ws.Range("E4").Value = "Holistic Fitness 6kg Rubber Slam Ball"
Set-TextCell "G4" "4,539"
Set-TextCell "H4" "60"
$ws.Range("I4").Value = "cod"
$ws.Range("J4").Value = "Delivered"

Set-TextCell "B5" "#FBZ-6C4S67Y"
$ws.Range("C5").Value = "Joyal K"
Set-TextCell "D5" "2025-03-04"
$ws.Range("E5").Value = "MEDIX Soft Medicine Ball (2), Rubber for Adults"
Set-TextCell "G5" "769"
Set-TextCell "H5" "60"
$ws.Range("I5").Value = "cod"
$ws.Range("J5").Value = "Shipped"

# ---- Insert the missing order as a new row 10 ----------------------------
# This pushes the blank separator row and the Summary block down by one.
$ws.Rows(10).Insert()

$ws.Range("A10").Value = 9
Set-TextCell "B10" "#FBZ-328BOQA"
$ws.Range("C10").Value = "Achu K"
Set-TextCell "D10" "2025-04-06"
$ws.Range("E10").Value = "MEDIX Soft Medicine Ball (2), Rubber for Adults"
$ws.Range("F10").Value = 1
Set-TextCell "G10" "629"
Set-TextCell "H10" "200"
$ws.Range("I10").Value = "cod"
$ws.Range("J10").Value = "Delivered"

# ---- Refresh the Summary totals (now at rows 12-15) ----------------------
$ws.Range("G12").Value = 9
$ws.Range("G13").Value = "₹28,125"
$ws.Range("G14").Value = "₹740"
$ws.Range("G15").Value = "₹27,385"
